$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old E column data (and old A column data we're about to overwrite)
$ws.Range("A1:E6").ClearContents()

# New header row
$ws.Range("A1").Value = "Test Case ID(s)"
$ws.Range("B1").Value = "Test Case Description"
$ws.Range("C1").Value = "Test Steps"

# Data row
$ws.Range("A2").Value = "TC_01"
$ws.Range("C2").Value = "1.Open URL ""http://localhost:8080/login"""
$ws.Range("C3").Value = "2.Type ""shana@gmail.com"" into ""//input[@placeholder='Enter your email']"""
$ws.Range("C4").Value = "3.Type ""shana"" into ""//input[@placeholder='Enter your password']"""
$ws.Range("C5").Value = "4.Click on the ""Sign In"" button ""//button[@type='submit']"""
$ws.Range("C6").Value = "5.Verify that the URL is ""http://localhost:8080/dashboard"""
$ws.Range("B2").Value = "Testing the Nesto App login page"

# Formatting: green fill on header row
$ws.Range("A1:C1").Interior.Color = 5287936

# Column widths (closest achievable values to the target stored widths of
# 13.44140625 / 35.44140625 / 81.6640625 characters)
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 34.666666666666664
$ws.Columns.Item(3).ColumnWidth = 80.83333333333334

# Selection matches target state
$ws.Range("C7").Select()
